$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume text columns keep their original "text" storage (many values
# look numeric, e.g. "4.021", and would otherwise be auto-converted to numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '27.687.82'
$ws.Range('E2').Value = '  -0.89%  '
$ws.Range('D3').Value = '1.739.80'
$ws.Range('E3').Value = '  -1.97%  '
$ws.Range('D4').Value = '1.007'
$ws.Range('E4').Value = '  +0.92%  '
$ws.Range('D5').Value = '330.62'
$ws.Range('E5').Value = '  -1.06%  '
$ws.Range('D6').Value = '1.003'
$ws.Range('E6').Value = '  +0.75%  '
$ws.Range('D7').Value = '0.3832'
$ws.Range('E7').Value = '  +0.80%  '
$ws.Range('D8').Value = '0.3340'
$ws.Range('E8').Value = '  -2.18%  '
$ws.Range('D9').Value = '45.28'
$ws.Range('E9').Value = '  -4.96%  '
$ws.Range('D10').Value = '1.095'
$ws.Range('E10').Value = '  -4.09%  '
$ws.Range('D11').Value = '0.07153'
$ws.Range('E11').Value = '  -3.64%  '
$ws.Range('D12').Value = '1.005'
$ws.Range('E12').Value = '  +0.86%  '
$ws.Range('D13').Value = '22.02'
$ws.Range('E13').Value = '  -3.61%  '
$ws.Range('D14').Value = '6.087'
$ws.Range('E14').Value = '  -4.49%  '
$ws.Range('D15').Value = '1.742.35'
$ws.Range('E15').Value = '  -1.51%  '
$ws.Range('D16').Value = '6.934'
$ws.Range('E16').Value = '  -2.51%  '
$ws.Range('D17').Value = '0.00001044'
$ws.Range('E17').Value = '  -3.25%  '
$ws.Range('D18').Value = '0.06559'
$ws.Range('E18').Value = '  -1.67%  '
$ws.Range('E19').Value = '  +0.59%  '
$ws.Range('D20').Value = '78.21'
$ws.Range('E20').Value = '  -5.07%  '
$ws.Range('D21').Value = '16.55'
$ws.Range('E21').Value = '  -5.00%  '
$ws.Range('D22').Value = '6.134'
$ws.Range('E22').Value = '  -4.46%  '
$ws.Range('D23').Value = '27.698.88'
$ws.Range('E23').Value = '  -0.72%  '
$ws.Range('D24').Value = '11.48'
$ws.Range('E24').Value = '  -5.01%  '
$ws.Range('D25').Value = '2.408'
$ws.Range('E25').Value = '  +1.46%  '
$ws.Range('D26').Value = '154.88'
$ws.Range('E26').Value = '  +0.60%  '
$ws.Range('D27').Value = '19.56'
$ws.Range('E27').Value = '  -5.78%  '
$ws.Range('D28').Value = '2.248'
$ws.Range('E28').Value = '  -7.55%  '
$ws.Range('D29').Value = '1.944.15'
$ws.Range('E29').Value = '  -1.29%  '
$ws.Range('D30').Value = '1.260'
$ws.Range('E30').Value = '  -12.34%  '
$ws.Range('D31').Value = '128.46'
$ws.Range('E31').Value = '  -4.36%  '
$ws.Range('D32').Value = '4.021'
$ws.Range('D33').Value = '5.716'
$ws.Range('E33').Value = '  -7.49%  '
$ws.Range('D34').Value = '0.08671'
$ws.Range('E34').Value = '  -1.32%  '
$ws.Range('D35').Value = '11.87'
$ws.Range('E35').Value = '  -7.36%  '
$ws.Range('D36').Value = '1.526'
$ws.Range('E36').Value = '  +0.63%  '
$ws.Range('B37').Value = 'TheSandbox'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D37').Value = '0.6419'
$ws.Range('E37').Value = '  -6.65%  '
$ws.Range('D38').Value = '0.02247'
$ws.Range('E38').Value = '  -7.75%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').Value = '5.064'
$ws.Range('E39').Value = '  -5.02%  '
$ws.Range('D40').Value = '0.06008'
$ws.Range('E40').Value = '  -5.29%  '
$ws.Range('D41').Value = '0.2074'
$ws.Range('E41').Value = '  -5.30%  '
$ws.Range('D42').Value = '1.185'
$ws.Range('E42').Value = '  -4.14%  '
$ws.Range('D43').Value = '1.003'
$ws.Range('E43').Value = '  +0.62%  '
$ws.Range('D44').Value = '7.884'
$ws.Range('E44').Value = '  -4.77%  '
$ws.Range('D45').Value = '13.62'
$ws.Range('E45').Value = '  -3.96%  '
$ws.Range('D46').Value = '3.790'
$ws.Range('E46').Value = '  -1.23%  '
$ws.Range('D47').Value = '0.5932'
$ws.Range('E47').Value = '  -5.88%  '
$ws.Range('D48').Value = '125.61'
$ws.Range('E48').Value = '  -4.82%  '
$ws.Range('D49').Value = '1.965'
$ws.Range('E49').Value = '  -5.82%  '
$ws.Range('D50').Value = '1.142'
$ws.Range('E50').Value = '  -0.75%  '
$ws.Range('D51').Value = '0.06925'
$ws.Range('E51').Value = '  -5.94%  '
